# Standardize "Multi KDMA" terminology in the definition column (column K/Z/AC/AG/AJ, row 2)
# of the Survey Delegation Variables sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = "X refers to the block number, Y refers to the DM number the participant saw. For Eval 8 and 9, participants saw 4 blocks with 3-4 DMs each (4 in the case of Multi KDMA). The following columns describe each page of the survey using this BX_DMY format."

$ws.Range("Z2").Value = "The name and alignment value of the fourth medic being compared in this comparison page - only applies to Multi KDMA"

$ws.Range("AC2").Value = "The response to the first forced choice question (aligned vs baseline, or follow the previous column for Multi KDMA)"

$ws.Range("AG2").Value = "The response to the second forced choice question (aligned vs misaligned, or follow the previous column for Multi KDMA)"

$ws.Range("AJ2").Value = "The alignment of the third two DMs being compared (Multi KDMA only)"
